# Story 71622: add a QA2 agency login link to cell A4 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New hyperlink cell: text + external link, matching the existing
# hyperlink cells already on the sheet (A2, B2, B3, B4, B5, B6).
$ws.Hyperlinks.Add($ws.Range("A4"), "https://login-agency-qa2.fnf.com/") | Out-Null

# Hyperlinks.Add() stamps the cell with an auto-generated hyperlink
# style; re-apply the workbook's named "Hyperlink" style explicitly so
# the cell reuses the same style index as the other hyperlink cells.
$ws.Range("A4").Style = "Hyperlink"

# Move the active selection to the newly populated cell.
$ws.Range("A4").Select() | Out-Null
